$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Note: Since BigOven search API*") {
        $p.Range.Font.Italic = $true
    }
}
